$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.861.21"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "2.236.57"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.97"
$ws.Range("E5").Value = "  -1.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.57"
$ws.Range("E6").Value = "  -3.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.567"
$ws.Range("E7").Value = "  -3.34%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  -6.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.98"
$ws.Range("E10").Value = "  -4.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0820"
$ws.Range("E11").Value = "  -2.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.32"
$ws.Range("E12").Value = "  -6.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("E13").Value = "  -3.09%  "
$ws.Range("D14").Value = "2.571.67"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.238.50"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.838"
$ws.Range("E16").Value = "  -4.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.04"
$ws.Range("E17").Value = "  -3.22%  "
$ws.Range("D18").Value = "43.719.89"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.77"
$ws.Range("E19").Value = "  -10.36%  "
$ws.Range("D20").Value = "0.0₃0960"
$ws.Range("E20").Value = "  -4.54%  "
$ws.Range("E21").Value = "  -4.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.83"
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.02"
$ws.Range("E23").Value = "  -4.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "233.27"
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.03"
$ws.Range("E25").Value = "  -8.39%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.20"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.18"
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.66"
$ws.Range("E29").Value = "  -3.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.98"
$ws.Range("E30").Value = "  -7.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.49"
$ws.Range("E31").Value = "  -2.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.90"
$ws.Range("E32").Value = "  -2.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0829"
$ws.Range("E33").Value = "  -6.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.69"
$ws.Range("E34").Value = "  -2.25%  "
$ws.Range("E35").Value = "  -3.23%  "
$ws.Range("E36").Value = "  +4.08%  "
$ws.Range("E37").Value = "  -4.10%  "
$ws.Range("E38").Value = "  -3.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.02"
$ws.Range("E39").Value = "  +2.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.61"
$ws.Range("E40").Value = "  -5.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.03"
$ws.Range("E41").Value = "  -9.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0308"
$ws.Range("E42").Value = "  -5.45%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "1.729.51"
$ws.Range("E44").Value = "  -4.57%  "
$ws.Range("E45").Value = "  -6.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "73.90"
$ws.Range("E46").Value = "  -2.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "80.26"
$ws.Range("E47").Value = "  -5.52%  "
$ws.Range("E48").Value = "  -5.12%  "
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.34"
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.39"
$ws.Range("E51").Value = "  -5.30%  "
